$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Status moved from "Ready for handoff" to "In Translation" everywhere it is shown
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# A handoff report was generated for the archive - record its name
$wsZhCn.Range("I2").Value = "TestHandoff_2016-12-05-10-09"
$wsDeDe.Range("I2").Value = "TestHandoff_2016-12-05-10-09"

# Let the columns that now hold different text re-fit their width
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsZhCn.Columns.Item(9).ColumnWidth = 28.1666666666667
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(9).ColumnWidth = 28.1666666666667
